$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix up cell styling BEFORE rewriting values, while we can still borrow
#    existing formats from the sheet (PasteSpecial -4122 = xlPasteFormats).
#    Target layout needs font19/border10 ("special" style) on columns C & G,
#    font19/border11 (same font, no left edge) on column D, and plain
#    bordered style everywhere else (incl. the new column N).
# ---------------------------------------------------------------------------

# C2/C3 and G2/G3 take on the "special" style currently sitting on E2/E3.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null

# D2/D3 take the same special style, then drop the left border edge.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Borders.Item(7).LineStyle = -4142
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Borders.Item(7).LineStyle = -4142

# E2/E3 revert to the plain bordered style (copy from B2/B3).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

# New column N: header cell gets the header style (copy from M1), data cells
# get the plain bordered style (copy from M2/M3).
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Rewrite header row with the new column order/labels.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value  = "Order Received Data and Time"
$ws.Cells.Item(1,2).Value  = "OrderID"
$ws.Cells.Item(1,3).Value  = "Emp ID-Order Assigned"
$ws.Cells.Item(1,4).Value  = "Assignee_QA"
$ws.Cells.Item(1,5).Value  = "Typist"
$ws.Cells.Item(1,6).Value  = "Typist QC"
$ws.Cells.Item(1,7).Value  = "Client"
$ws.Cells.Item(1,8).Value  = "Lob"
$ws.Cells.Item(1,9).Value  = "Process"
$ws.Cells.Item(1,10).Value = "Product Name"
$ws.Cells.Item(1,11).Value = "State"
$ws.Cells.Item(1,12).Value = "County"
$ws.Cells.Item(1,13).Value = "Status"
$ws.Cells.Item(1,14).Value = "Tier"

# ---------------------------------------------------------------------------
# 3. Rewrite data rows 2 & 3.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value  = 45436.041666666664
$ws.Cells.Item(2,2).Value  = "RSS18-001"
$ws.Cells.Item(2,3).Value  = "SIPL5316"
$ws.Cells.Item(2,4).Value  = "SIPL5688"
$ws.Cells.Item(2,5).Value  = "SIPL5317"
$ws.Cells.Item(2,6).Value  = "SIPL5317"
$ws.Cells.Item(2,7).Value  = "Radian Settlement Services Inc"
$ws.Cells.Item(2,8).Value  = "Title"
$ws.Cells.Item(2,9).Value  = "Search & Typing"
$ws.Cells.Item(2,10).Value = "Radian Search and Type"
$ws.Cells.Item(2,11).Value = "AL"
$ws.Cells.Item(2,12).Value = "Shelby"
$ws.Cells.Item(2,13).Value = "WIP"
$ws.Cells.Item(2,14).Value = ""

$ws.Cells.Item(3,1).Value  = 45439.083333333336
$ws.Cells.Item(3,2).Value  = "RSS18-002"
$ws.Cells.Item(3,3).Value  = "SIPL5316"
$ws.Cells.Item(3,4).Value  = "SIPL5688"
$ws.Cells.Item(3,5).Value  = "SIPL5317"
$ws.Cells.Item(3,6).Value  = "SIPL5317"
$ws.Cells.Item(3,7).Value  = "Radian Settlement Services Inc"
$ws.Cells.Item(3,8).Value  = "Title"
$ws.Cells.Item(3,9).Value  = "Search & Typing"
$ws.Cells.Item(3,10).Value = "Update Search"
$ws.Cells.Item(3,11).Value = "FL"
$ws.Cells.Item(3,12).Value = "Clay"
$ws.Cells.Item(3,13).Value = "WIP"
$ws.Cells.Item(3,14).Value = ""

# ---------------------------------------------------------------------------
# 4. Column widths to match the new layout. The host engine re-quantises
#    ColumnWidth to the nearest 1/6 of a character, so the inputs below are
#    pre-compensated (target - 5/6, rounded to the nearest 1/6) so the saved
#    width lands as close as possible to the real template's values
#    (20.5546875 / 8.33203125 / 26.21875 / 10 / 16 / 20.44140625).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(5).ColumnWidth = 7.500000000000001
$ws.Columns.Item(7).ColumnWidth = 25.333333333333336
$ws.Columns.Item(8).ColumnWidth = 9.166666666666666
$ws.Columns.Item(9).ColumnWidth = 15.166666666666666
$ws.Columns.Item(10).ColumnWidth = 19.666666666666668

# ---------------------------------------------------------------------------
# 5. Selection, matching the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("I5").Select() | Out-Null
